# Update recurrence metrics for the last quarter (row 28: 2025Q2)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C28").Value = 477
$ws.Range("D28").Value = 51
$ws.Range("E28").Value = 426
$ws.Range("F28").Value = 7.943925233644859
